$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($sheet, $addr, $val) {
    # Force the cell to be written as text so values such as "1.00",
    # "0.512" or "212.84" are not auto-converted into numbers (which would
    # change their stored representation), matching the original file's
    # inline-string cells. The temporary "@" number format is removed
    # again right away so the cell keeps its original (default) style.
    $sheet.Range($addr).NumberFormat = "@"
    $sheet.Range($addr).Value = $val
    $sheet.Range($addr).Style = "Normal"
}

# --- Simple per-row updates (Price in column D, Volume(1h) in column E) ---

Set-TextValue $ws "D2" "26.257.89"
$ws.Range("E2").Value = "  +3.64%  "

Set-TextValue $ws "D3" "1.606.67"
$ws.Range("E3").Value = "  +2.14%  "

$ws.Range("E4").Value = "  -0.59%  "

Set-TextValue $ws "D5" "212.84"
$ws.Range("E5").Value = "  +2.26%  "

$ws.Range("E6").Value = "  -0.68%  "

$ws.Range("E7").Value = "  +1.55%  "

$ws.Range("E8").Value = "  +2.09%  "

$ws.Range("E9").Value = "  +1.64%  "

Set-TextValue $ws "D10" "18.09"
$ws.Range("E10").Value = "  +0.78%  "

Set-TextValue $ws "D11" "0.0815"
$ws.Range("E11").Value = "  +4.28%  "

Set-TextValue $ws "D12" "1.830.12"
$ws.Range("E12").Value = "  +2.15%  "

Set-TextValue $ws "D13" "1.606.26"
$ws.Range("E13").Value = "  +1.92%  "

$ws.Range("E14").Value = "  -0.41%  "

Set-TextValue $ws "D15" "0.512"
$ws.Range("E15").Value = "  +1.06%  "

Set-TextValue $ws "D16" "26.203.75"
$ws.Range("E16").Value = "  +3.44%  "

Set-TextValue $ws "D17" "60.62"
$ws.Range("E17").Value = "  +1.37%  "

$sub3 = [char]0x2083
Set-TextValue $ws "D18" "0.0${sub3}0727"
$ws.Range("E18").Value = "  +2.12%  "

$ws.Range("E19").Value = "  -0.56%  "

Set-TextValue $ws "D20" "198.52"
$ws.Range("E20").Value = "  +6.80%  "

$ws.Range("E21").Value = "  +2.31%  "

Set-TextValue $ws "D22" "9.38"
$ws.Range("E22").Value = "  +0.52%  "

$ws.Range("E23").Value = "  +1.80%  "

Set-TextValue $ws "D24" "142.34"
$ws.Range("E24").Value = "  +0.75%  "

$ws.Range("E25").Value = "  +3.39%  "

# --- Row 26 / Row 27: contents swap (Stellar <-> BinanceUSD), plus value updates ---

$ws.Range("B26").Value = "Stellar"
$ws.Range("C26").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws "D26" "0.127"
$ws.Range("E26").Value = "  -0.91%  "

$ws.Range("B27").Value = "BinanceUSD"
$ws.Range("C27").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue $ws "D27" "1.00"
$ws.Range("E27").Value = "  -0.58%  "

# --- Continue simple updates ---

Set-TextValue $ws "D28" "15.17"
$ws.Range("E28").Value = "  +1.86%  "

Set-TextValue $ws "D29" "6.47"
$ws.Range("E29").Value = "  +0.07%  "

Set-TextValue $ws "D30" "1.17"
$ws.Range("E30").Value = "  +1.36%  "

$ws.Range("E31").Value = "  +2.10%  "

$ws.Range("E32").Value = "  +2.71%  "

$ws.Range("E33").Value = "  -0.17%  "

$ws.Range("E34").Value = "  +0.75%  "

$ws.Range("E35").Value = "  +4.73%  "

Set-TextValue $ws "D36" "1.108.63"
$ws.Range("E36").Value = "  +2.07%  "

$ws.Range("E37").Value = "  -0.53%  "

$ws.Range("E38").Value = "  +1.00%  "

$ws.Range("E39").Value = "  +0.39%  "

$ws.Range("E40").Value = "  +0.68%  "

Set-TextValue $ws "D41" "0.500"
$ws.Range("E41").Value = "  +1.02%  "

Set-TextValue $ws "D42" "0.776"
$ws.Range("E42").Value = "  +3.42%  "

Set-TextValue $ws "D43" "1.741.95"
$ws.Range("E43").Value = "  +2.14%  "

$ws.Range("E44").Value = "  +0.91%  "

Set-TextValue $ws "D45" "92.42"
$ws.Range("E45").Value = "  -0.67%  "

$ws.Range("E46").Value = "  +0.69%  "

Set-TextValue $ws "D47" "1.56"
$ws.Range("E47").Value = "  +9.93%  "

Set-TextValue $ws "D48" "53.48"
$ws.Range("E48").Value = "  +1.14%  "

Set-TextValue $ws "D50" "0.409"
$ws.Range("E50").Value = "  +0.37%  "

$ws.Range("E51").Value = "  -0.33%  "
